$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "70.516.24"
$ws.Range("E2").Value = "  +0.74%  "

# Row 3
Set-TextValue "D3" "3.768.26"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
Set-TextValue "D5" "621.04"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
Set-TextValue "D6" "181.30"
$ws.Range("E6").Value = "  +2.83%  "

# Row 7
Set-TextValue "D7" "3.766.00"
$ws.Range("E7").Value = "  -0.29%  "

# Row 8
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  -1.88%  "

# Row 10
Set-TextValue "D10" "0.170"
$ws.Range("E10").Value = "  +1.69%  "

# Row 11
Set-TextValue "D11" "6.39"
$ws.Range("E11").Value = "  +0.37%  "

# Row 12
Set-TextValue "D12" "0.484"
$ws.Range("E12").Value = "  -3.21%  "

# Row 13
Set-TextValue "D13" "40.50"
$ws.Range("E13").Value = "  +0.21%  "

# Row 14
Set-TextValue "D14" "0.0000259"
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
Set-TextValue "D15" "4.393.07"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
Set-TextValue "D16" "3.769.99"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17
Set-TextValue "D17" "70.717.77"
$ws.Range("E17").Value = "  +0.64%  "

# Row 18
Set-TextValue "D18" "7.62"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19
$ws.Range("E19").Value = "  -1.84%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "16.60"
$ws.Range("E20").Value = "  -0.49%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "507.44"
$ws.Range("E21").Value = "  -3.14%  "

# Row 22
Set-TextValue "D22" "9.26"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
Set-TextValue "D23" "0.725"
$ws.Range("E23").Value = "  -2.41%  "

# Row 24
$ws.Range("E24").Value = "  +4.09%  "

# Row 25
Set-TextValue "D25" "87.34"
$ws.Range("E25").Value = "  -1.37%  "

# Row 26
$ws.Range("E26").Value = "  -2.97%  "

# Row 27
Set-TextValue "D27" "11.27"
$ws.Range("E27").Value = "  +3.43%  "

# Row 28
Set-TextValue "D28" "0.0000139"
$ws.Range("E28").Value = "  +12.67%  "

# Row 29
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
$ws.Range("E30").Value = "  -0.46%  "

# Row 31
$ws.Range("E31").Value = "  +2.69%  "

# Row 32
Set-TextValue "D32" "7.88"
$ws.Range("E32").Value = "  -0.20%  "

# Row 33
Set-TextValue "D33" "30.72"
$ws.Range("E33").Value = "  -3.81%  "

# Row 34
$ws.Range("E34").Value = "  +0.23%  "

# Row 35
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
Set-TextValue "D36" "1.06"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37
Set-TextValue "D37" "6.14"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
Set-TextValue "D38" "0.352"
$ws.Range("E38").Value = "  +2.81%  "

# Row 39
$ws.Range("E39").Value = "  +4.98%  "

# Row 40
$ws.Range("E40").Value = "  +18.47%  "

# Row 41
Set-TextValue "D41" "2.09"
$ws.Range("E41").Value = "  -3.17%  "

# Row 42
Set-TextValue "D42" "50.01"
$ws.Range("E42").Value = "  -2.88%  "

# Row 43
Set-TextValue "D43" "434.92"
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
Set-TextValue "D44" "44.35"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
Set-TextValue "D45" "8.65"
$ws.Range("E45").Value = "  -2.05%  "

# Row 46
Set-TextValue "D46" "2.974.32"
$ws.Range("E46").Value = "  -5.17%  "

# Row 47
Set-TextValue "D47" "0.0365"
$ws.Range("E47").Value = "  -0.61%  "

# Row 48
Set-TextValue "D48" "27.36"
$ws.Range("E48").Value = "  -1.13%  "

# Row 49
$ws.Range("E49").Value = "  -0.03%  "

# Row 50
Set-TextValue "D50" "136.74"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
Set-TextValue "D51" "2.50"
$ws.Range("E51").Value = "  -0.85%  "

Write-Host "All updates applied"